# Apply "Penalty Reward System" update to the forecast workbook.
# - "Forecast Comparison" sheet: shift each week's date forward by one week
#   and refresh the MyForecast (column D) figures.
# - "Summary" sheet: refresh the derived summary statistics to match.
#
# The text-looking values in column B (dates, "NNN units", plain numbers
# stored as text, etc.) must stay plain text cells, exactly as authored in
# the source file. Assigning a date-like / numeric-like string straight to
# `.Value` makes Excel's input parser silently re-type it (date serial /
# number), so instead we briefly turn the literal into a `="..."` text
# formula (never auto-converted) and immediately collapse it back down to a
# constant with a values-only copy/paste - exactly what "Paste Special ->
# Values" does when cleaning up a formula by hand, and it leaves the cell's
# original (default) formatting completely untouched.

function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Value = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: columns B (Week_Start_Date) and D (MyForecast) ---
# row -> @(newDate, newMyForecast)
$rows = @{
    2  = @("2025-01-12", 63)
    3  = @("2025-01-19", 63)
    4  = @("2025-01-26", 64)
    5  = @("2025-02-02", 66)
    6  = @("2025-02-09", 69)
    7  = @("2025-02-16", 73)
    8  = @("2025-02-23", 51)
    9  = @("2025-03-02", 78)
    10 = @("2025-03-09", 72)
    11 = @("2025-03-16", 67)
    12 = @("2025-03-23", 65)
    13 = @("2025-03-30", 68)
    14 = @("2025-04-06", 73)
    15 = @("2025-04-13", 49)
    16 = @("2025-04-20", 71)
    17 = @("2025-04-27", 47)
}

foreach ($r in ($rows.Keys | Sort-Object)) {
    $vals = $rows[$r]
    Set-TextValue $wsForecast.Cells.Item($r, 2) $vals[0]
    $wsForecast.Cells.Item($r, 4).Value = $vals[1]
}

# --- Summary sheet updates (all column-B cells here are stored as text) ---
$summaryUpdates = [ordered]@{
    "B2"  = "2023-01-08 to 2025-01-05"
    "B8"  = "6160 units"
    "B9"  = "1040"
    "B10" = "527"
    "B11" = "256"
    "B12" = "78"
    "B13" = "2025-03-02"
    "B14" = "47"
    "B15" = "2025-04-27"
}

foreach ($addr in $summaryUpdates.Keys) {
    Set-TextValue $wsSummary.Range($addr) $summaryUpdates[$addr]
}
